$d = $word.ActiveDocument

# --- Change 1: split the "humanized_is_absent_info_about_legal_proceedings..." field
#     into three runs with a new bookmark wrapped around the middle part.
$rng1 = $d.Content
$found1 = $rng1.Find.Execute(
    "{issue.humanized_is_absent_info_about_legal_proceedings_as_defendant_for_more_than_30_pct_of_net_assets}",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found1) {
    $start1 = $rng1.Start
    $end1 = $rng1.End

    $prefix = "{issue."
    $suffix = "}"

    $midStart = $start1 + $prefix.Length
    $midEnd = $end1 - $suffix.Length

    $midRange1 = $d.Range($midStart, $midEnd)
    $d.Bookmarks.Add("__DdeLink__20254_1008482545", $midRange1) | Out-Null
}

# --- Change 2: replace the "Положительное/ отрицательное " placeholder text
#     with the final DOMC conclusion merge field.
$rng2 = $d.Content
$rng2.Find.Execute(
    "Положительное/ отрицательное ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "{issue.humanized_final_documents_operations_management_conclusion}", 2) | Out-Null
